# Scheduled-runner style refresh of market-price-derived columns
# (currentAveragePrice*, Leve*Price*, Leve*Profit*) across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2635.7334
$ws.Range("I33").Value = 3004
$ws.Range("K33").Value = 3004
$ws.Range("M33").Value = -2775

$ws.Range("H43").Value = 143877.28
$ws.Range("J43").Value = 1338.2
$ws.Range("L43").Value = 1338.2
$ws.Range("N43").Value = -1476.2

$ws.Range("H80").Value = 3049
$ws.Range("I80").Value = 3742
$ws.Range("J80").Value = 1576.375
$ws.Range("K80").Value = 11226
$ws.Range("L80").Value = 4729.125
$ws.Range("M80").Value = -10228
$ws.Range("N80").Value = -6725.125

$ws.Range("H83").Value = 3049
$ws.Range("I83").Value = 3742
$ws.Range("J83").Value = 1576.375
$ws.Range("K83").Value = 33678
$ws.Range("L83").Value = 14187.375
$ws.Range("M83").Value = -28686
$ws.Range("N83").Value = -24171.375

$ws.Range("H88").Value = 18406.16
$ws.Range("I88").Value = 1502.75
$ws.Range("J88").Value = 26360.705
$ws.Range("K88").Value = 1502.75
$ws.Range("L88").Value = 26360.705
$ws.Range("M88").Value = -1096.75
$ws.Range("N88").Value = -27172.705

$ws.Range("H91").Value = 18406.16
$ws.Range("I91").Value = 1502.75
$ws.Range("J91").Value = 26360.705
$ws.Range("K91").Value = 1502.75
$ws.Range("L91").Value = 26360.705
$ws.Range("M91").Value = -98.75
$ws.Range("N91").Value = -29168.705

$ws.Range("H107").Value = 1380.9333
$ws.Range("I107").Value = 1382.8
$ws.Range("J107").Value = 1377.2
$ws.Range("K107").Value = 1382.8
$ws.Range("L107").Value = 1377.2
$ws.Range("M107").Value = 537.2
$ws.Range("N107").Value = -5217.2

$ws.Range("H132").Value = 1881.75
$ws.Range("I132").Value = 1798.826
$ws.Range("J132").Value = 2517.5
$ws.Range("K132").Value = 5396.478
$ws.Range("L132").Value = 7552.5
$ws.Range("M132").Value = -2866.478
$ws.Range("N132").Value = -12612.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1335.8292
$ws.Range("I2").Value = 1236.8235
$ws.Range("K2").Value = 1236.8235
$ws.Range("M2").Value = -1123.8235

$ws.Range("H45").Value = 798.2222
$ws.Range("I45").Value = 713.125
$ws.Range("K45").Value = 713.125
$ws.Range("M45").Value = -336.125

$ws.Range("H61").Value = 4773.346
$ws.Range("I61").Value = 3008.818
$ws.Range("J61").Value = 7838.0527
$ws.Range("K61").Value = 3008.818
$ws.Range("L61").Value = 7838.0527
$ws.Range("M61").Value = -2796.818
$ws.Range("N61").Value = -8262.0527

$ws.Range("H63").Value = 2985.5925
$ws.Range("I63").Value = 2322.7273
$ws.Range("K63").Value = 2322.7273
$ws.Range("M63").Value = -1636.7273

$ws.Range("H66").Value = 2985.5925
$ws.Range("I66").Value = 2322.7273
$ws.Range("K66").Value = 11613.6365
$ws.Range("M66").Value = -8181.636500000001

$ws.Range("H88").Value = 742.94446
$ws.Range("I88").Value = 630.75
$ws.Range("J88").Value = 967.3333
$ws.Range("K88").Value = 630.75
$ws.Range("L88").Value = 967.3333
$ws.Range("M88").Value = -224.75
$ws.Range("N88").Value = -1779.3333

$ws.Range("H91").Value = 742.94446
$ws.Range("I91").Value = 630.75
$ws.Range("J91").Value = 967.3333
$ws.Range("K91").Value = 630.75
$ws.Range("L91").Value = 967.3333
$ws.Range("M91").Value = 773.25
$ws.Range("N91").Value = -3775.3333

$ws.Range("H116").Value = 1335.8292
$ws.Range("I116").Value = 1236.8235
$ws.Range("K116").Value = 1236.8235
$ws.Range("M116").Value = 1057.1765

$ws.Range("H122").Value = 45771.617
$ws.Range("I122").Value = 3015.5386
$ws.Range("J122").Value = 115250.25
$ws.Range("K122").Value = 9046.6158
$ws.Range("L122").Value = 345750.75
$ws.Range("M122").Value = -6596.6158
$ws.Range("N122").Value = -350650.75

$ws.Range("H136").Value = 4773.346
$ws.Range("I136").Value = 3008.818
$ws.Range("J136").Value = 7838.0527
$ws.Range("K136").Value = 9026.454000000002
$ws.Range("L136").Value = 23514.1581
$ws.Range("M136").Value = -6476.454000000002
$ws.Range("N136").Value = -28614.1581

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1335.8292
$ws.Range("I3").Value = 1236.8235
$ws.Range("K3").Value = 1236.8235
$ws.Range("M3").Value = -1122.8235

$ws.Range("H107").Value = 14801.1
$ws.Range("I107").Value = 17377.484
$ws.Range("K107").Value = 17377.484
$ws.Range("M107").Value = -15457.484

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 30333
$ws.Range("I42").Value = 20499.5
$ws.Range("K42").Value = 20499.5
$ws.Range("M42").Value = -19906.5

$ws.Range("H86").Value = 2445.6365
$ws.Range("I86").Value = 2567.5334
$ws.Range("J86").Value = 2184.4285
$ws.Range("K86").Value = 2567.5334
$ws.Range("L86").Value = 2184.4285
$ws.Range("M86").Value = -1444.5334
$ws.Range("N86").Value = -4430.4285

$ws.Range("H89").Value = 2445.6365
$ws.Range("I89").Value = 2567.5334
$ws.Range("J89").Value = 2184.4285
$ws.Range("K89").Value = 12837.667
$ws.Range("L89").Value = 10922.1425
$ws.Range("M89").Value = -7221.666999999999
$ws.Range("N89").Value = -22154.1425

$ws.Range("H99").Value = 5348.1953
$ws.Range("I99").Value = 4853.5938
$ws.Range("J99").Value = 7106.778
$ws.Range("K99").Value = 4853.5938
$ws.Range("L99").Value = 7106.778
$ws.Range("M99").Value = -3355.5938
$ws.Range("N99").Value = -10102.778

$ws.Range("H126").Value = 5348.1953
$ws.Range("I126").Value = 4853.5938
$ws.Range("J126").Value = 7106.778
$ws.Range("K126").Value = 14560.7814
$ws.Range("L126").Value = 21320.334
$ws.Range("M126").Value = -12090.7814
$ws.Range("N126").Value = -26260.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9426.154
$ws.Range("I80").Value = 4089.9092
$ws.Range("J80").Value = 13339.4
$ws.Range("K80").Value = 4089.9092
$ws.Range("L80").Value = 13339.4
$ws.Range("M80").Value = -3091.9092
$ws.Range("N80").Value = -15335.4

$ws.Range("H83").Value = 9426.154
$ws.Range("I83").Value = 4089.9092
$ws.Range("J83").Value = 13339.4
$ws.Range("K83").Value = 20449.546
$ws.Range("L83").Value = 66697
$ws.Range("M83").Value = -15457.546
$ws.Range("N83").Value = -76681

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 57582.5
$ws.Range("J38").Value = 57582.5
$ws.Range("L38").Value = 57582.5
$ws.Range("N38").Value = -58402.5

$ws.Range("H41").Value = 9000
$ws.Range("I41").Value = 9000
$ws.Range("K41").Value = 9000
$ws.Range("M41").Value = -8562

$ws.Range("H46").Value = 1172.25
$ws.Range("I46").Value = 1229.6666
$ws.Range("K46").Value = 1229.6666
$ws.Range("M46").Value = -1041.6666

$ws.Range("H68").Value = 3175.5293
$ws.Range("J68").Value = 4132.8335
$ws.Range("L68").Value = 4132.8335
$ws.Range("N68").Value = -5630.8335

$ws.Range("H71").Value = 3175.5293
$ws.Range("J71").Value = 4132.8335
$ws.Range("L71").Value = 20664.1675
$ws.Range("N71").Value = -28152.1675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 23445
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 23445
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 23445
$ws.Range("N18").Value = -23791
$ws.Range("M18").ClearContents()

$ws.Range("H62").Value = 147934.86
$ws.Range("J62").Value = 4332.3335
$ws.Range("L62").Value = 4332.3335
$ws.Range("N62").Value = -5580.3335

$ws.Range("H65").Value = 147934.86
$ws.Range("J65").Value = 4332.3335
$ws.Range("L65").Value = 21661.6675
$ws.Range("N65").Value = -27901.6675

$ws.Range("H126").Value = 22819.9
$ws.Range("I126").Value = 27513.812
$ws.Range("K126").Value = 82541.436
$ws.Range("M126").Value = -80071.436
